$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell reference, new text value, whether the value must be forced
# to Text format (so Excel doesn't reinterpret a numeric-looking string like
# "236.98" or "1.000" as a number and strip / reformat it).
$updates = @(
    @('D2', '29.094.08', 0),
    @('E2', '  -2.79%  ', 0),
    @('D3', '1.843.65', 0),
    @('E3', '  -1.76%  ', 0),
    @('E4', '  +0.07%  ', 0),
    @('D5', '0.6963', 1),
    @('E5', '  -6.23%  ', 0),
    @('D6', '236.98', 1),
    @('E6', '  -2.28%  ', 0),
    @('E7', '  +0.06%  ', 0),
    @('D8', '0.3022', 1),
    @('E8', '  -4.06%  ', 0),
    @('D9', '0.07395', 1),
    @('E9', '  +2.44%  ', 0),
    @('D10', '23.25', 1),
    @('E10', '  -6.19%  ', 0),
    @('D11', '0.08102', 1),
    @('E11', '  -2.95%  ', 0),
    @('B12', 'Polygon', 0),
    @('C12', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', 0),
    @('D12', '0.7213', 1),
    @('E12', '  -3.95%  ', 0),
    @('B13', 'WrappedEther', 0),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', 0),
    @('D13', '1.827.83', 0),
    @('E13', '  -2.14%  ', 0),
    @('D14', '5.171', 1),
    @('E14', '  -4.17%  ', 0),
    @('D15', '88.80', 1),
    @('E15', '  -3.78%  ', 0),
    @('D16', '29.157.60', 0),
    @('E16', '  -2.58%  ', 0),
    @('D17', '5.764', 1),
    @('E17', '  -5.93%  ', 0),
    @('D18', '240.34', 1),
    @('E18', '  -3.09%  ', 0),
    @('D19', '0.000007644', 1),
    @('E19', '  -2.62%  ', 0),
    @('D20', '12.95', 1),
    @('E20', '  -4.53%  ', 0),
    @('D21', '1.000', 1),
    @('E21', '  +0.03%  ', 0),
    @('B22', 'WrappedliquidstakedEther2.0', 0),
    @('C22', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', 0),
    @('D22', '2.102.16', 0),
    @('E22', '  -1.41%  ', 0),
    @('B23', 'BinanceUSD', 0),
    @('C23', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', 0),
    @('D23', '1.001', 1),
    @('E23', '  +0.11%  ', 0),
    @('D24', '7.598', 1),
    @('E24', '  -5.06%  ', 0),
    @('D25', '0.1464', 1),
    @('E25', '  -5.81%  ', 0),
    @('D26', '161.36', 1),
    @('E26', '  -2.76%  ', 0),
    @('D27', '8.967', 1),
    @('E27', '  -3.52%  ', 0),
    @('D28', '17.98', 1),
    @('E28', '  -3.72%  ', 0),
    @('D29', '1.922', 1),
    @('E29', '  -5.16%  ', 0),
    @('D30', '1.377', 1),
    @('E30', '  -7.10%  ', 0),
    @('D31', '4.436', 1),
    @('E31', '  -3.33%  ', 0),
    @('D32', '1.487', 1),
    @('E32', '  -3.22%  ', 0),
    @('D33', '3.996', 1),
    @('E33', '  -5.62%  ', 0),
    @('D34', '0.05184', 1),
    @('E34', '  -3.34%  ', 0),
    @('D35', '1.179', 1),
    @('E35', '  -5.15%  ', 0),
    @('D36', '0.7095', 1),
    @('E36', '  -5.47%  ', 0),
    @('D37', '0.9991', 1),
    @('E37', '  -0.23%  ', 0),
    @('D38', '2.649', 1),
    @('E38', '  -2.11%  ', 0),
    @('D39', '0.01869', 1),
    @('E39', '  -4.89%  ', 0),
    @('D40', '2.669', 1),
    @('E40', '  -3.08%  ', 0),
    @('D41', '0.9111', 1),
    @('E41', '  +5.63%  ', 0),
    @('D42', '0.4268', 1),
    @('E42', '  -6.23%  ', 0),
    @('D43', '5.887', 1),
    @('E43', '  -4.02%  ', 0),
    @('D44', '1.057.90', 0),
    @('E44', '  -5.75%  ', 0),
    @('D45', '69.59', 1),
    @('E45', '  -4.13%  ', 0),
    @('D46', '0.9999', 1),
    @('E46', '  -0.05%  ', 0),
    @('D47', '101.66', 1),
    @('E47', '  -2.64%  ', 0),
    @('B48', 'RocketPoolETH', 0),
    @('C48', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', 0),
    @('D48', '2.009.42', 0),
    @('E48', '  -0.79%  ', 0),
    @('B49', 'RenderToken', 0),
    @('C49', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', 0),
    @('D49', '1.741', 1),
    @('E49', '  -6.83%  ', 0),
    @('D50', '9.252', 1),
    @('E50', '  -2.71%  ', 0),
    @('B51', 'Aptos', 0),
    @('C51', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', 0),
    @('D51', '7.097', 1),
    @('E51', '  -6.87%  ', 0)
)

foreach ($u in $updates) {
    $cellRef   = $u[0]
    $newValue  = $u[1]
    $forceText = $u[2]

    $rng = $ws.Range($cellRef)
    if ($forceText -eq 1) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $newValue
}
